$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# --- Apply formatting to the new Cannister A / Cannister B data block by
# --- copying the existing per-column formats used elsewhere on the sheet, so
# --- the saved file reuses the same style indices Excel itself would use.

# Column B (index numbers): General / Arial 10 style, same as A2/B2.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B20:B33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column C (Cannister A results): 0% / Arial 10 style, same as C2.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C20:C33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column D (Cannister B results): 0% / Calibri "Percent" style, same as F2 --
# but only for the cells that hold a number; the "NA" cells stay unformatted.
# (pasted one contiguous area at a time; multi-area paste only applies to the
# first area under this runtime)
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D24:D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D30:D33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 22's "NA" text cell in column C uses the plain General/Arial style
# (same as column B), not the percentage style used by the other C cells.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Table headers (row 19) ---
$ws.Range("C19").Value2 = "Cannister A"
$ws.Range("D19").Value2 = "Cannister B"

# --- Data rows ---
$ws.Range("B20").Value2 = 1
$ws.Range("C20").Value2 = 0.94
$ws.Range("D20").Value2 = 0.9

$ws.Range("B21").Value2 = 2
$ws.Range("C21").Value2 = 0.92
$ws.Range("D21").Value2 = "NA"

$ws.Range("B22").Value2 = 3
$ws.Range("C22").Value2 = "NA"
$ws.Range("D22").Value2 = 0.9

$ws.Range("B23").Value2 = 4
$ws.Range("C23").Value2 = 0.88
$ws.Range("D23").Value2 = "NA"

$ws.Range("B24").Value2 = 5
$ws.Range("C24").Value2 = 0.87
$ws.Range("D24").Value2 = 0.87

$ws.Range("B25").Value2 = 6
$ws.Range("C25").Value2 = 0.92
$ws.Range("D25").Value2 = 0.89

$ws.Range("B26").Value2 = 7
$ws.Range("C26").Value2 = 0.96
$ws.Range("D26").Value2 = 0.96

$ws.Range("B27").Value2 = 7
$ws.Range("C27").Value2 = 0.94
$ws.Range("D27").Value2 = 0.97

$ws.Range("B28").Value2 = 8
$ws.Range("C28").Value2 = 1
$ws.Range("D28").Value2 = 1

$ws.Range("B29").Value2 = 9
$ws.Range("C29").Value2 = 1
$ws.Range("D29").Value2 = "NA "

$ws.Range("B30").Value2 = 10
$ws.Range("C30").Value2 = 1
$ws.Range("D30").Value2 = 1

$ws.Range("B31").Value2 = 11
$ws.Range("C31").Value2 = 1
$ws.Range("D31").Value2 = 1

$ws.Range("B32").Value2 = 12
$ws.Range("C32").Value2 = 1
$ws.Range("D32").Value2 = 1

$ws.Range("B33").Value2 = 13
$ws.Range("C33").Value2 = 0.99
$ws.Range("D33").Value2 = 0.99

# --- Selection, matching the author's final cursor position ---
$ws.Range("C29").Select() | Out-Null
